# Apply categorization fix: reclassify two "Other" transactions as "Income"
# and update the Summary sheet to reflect the new totals.

$wb = $excel.ActiveWorkbook

# --- Transactions sheet ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("D7").Value = "Income"
$wsTrans.Range("D8").Value = "Income"

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")

# Update Income total (2000 -> 1805)
$wsSummary.Range("B3").Value = 1805

# "Other" category row is gone; Shopping/Transport shift up one row.
$wsSummary.Range("A4").Value = "Shopping"
$wsSummary.Range("B4").Value = -45
$wsSummary.Range("A5").Value = "Transport"
$wsSummary.Range("B5").Value = -22.75

# Remove the now-empty trailing row 6 (previously Transport) entirely,
# shifting rows up so the used range shrinks to A1:B5.
$wsSummary.Rows.Item(6).Delete()
